# Add all updated files and figures for 2021 reconstruction
#
# This adds a new row 60 of data to each of the four worksheets
# (Escapement, Total Catch, Run Size, Run Size no Offshore) and
# updates some of the existing row 59 values on the latter three
# sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: Escapement ---
$ws1 = $wb.Worksheets.Item("Escapement")
$ws1.Range("A60").Value = 4688089.278375
$ws1.Range("B60").Value = 3230621.48119301
$ws1.Range("C60").Value = 2780528.141222
$ws1.Range("D60").Value = 1849311.35638
$ws1.Range("E60").Value = 2902367.4527424

# --- Sheet 2: Total Catch ---
$ws2 = $wb.Worksheets.Item("Total Catch")
$ws2.Range("A59").Value = 5835696.59801433
$ws2.Range("B59").Value = 2174689.96366779
$ws2.Range("C59").Value = 5990051.95581967
$ws2.Range("D59").Value = 13207490.5318939
$ws2.Range("E59").Value = 3332673.97612231

$ws2.Range("A60").Value = 4173540.65451519
$ws2.Range("B60").Value = 2824542.76335536
$ws2.Range("C60").Value = 4315536.2001339
$ws2.Range("D60").Value = 5885686.24853126
$ws2.Range("E60").Value = 7805217.15984679

# --- Sheet 3: Run Size ---
$ws3 = $wb.Worksheets.Item("Run Size")
$ws3.Range("A59").Value = 9866664.59768534
$ws3.Range("B59").Value = 4561207.9636928
$ws3.Range("C59").Value = 10102211.9563297
$ws3.Range("D59").Value = 15597218.5321769
$ws3.Range("E59").Value = 5078613.97604791

$ws3.Range("A60").Value = 8877060.6548902
$ws3.Range("B60").Value = 6061446.76354836
$ws3.Range("C60").Value = 7112070.1993559
$ws3.Range("D60").Value = 7717882.24891126
$ws3.Range("E60").Value = 10665147.1595892

# --- Sheet 4: Run Size no Offshore ---
$ws4 = $wb.Worksheets.Item("Run Size no Offshore")
$ws4.Range("A59").Value = 9813005.95416345
$ws4.Range("B59").Value = 4535465.94022568
$ws4.Range("C59").Value = 10047839.8769871

$ws4.Range("A60").Value = 8453208.9993308
$ws4.Range("B60").Value = 5768722.89931442
$ws4.Range("C60").Value = 6801010.79850064
$ws4.Range("D60").Value = 7315452.9239985
$ws4.Range("E60").Value = 10108062.8180798
